$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New crypto price/volume snapshot data (rows 2-51, columns B:E)
# Row 8 now holds a newly-added coin (LidoStakedEther), shifting the
# previously lower-ranked coins down by one row; Cronos drops off the list.
$data = @(
    @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.807.52', '  -0.54%  '),
    @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.895.66', '  +0.15%  '),
    @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  +0.09%  '),
    @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.7644', '  +3.85%  '),
    @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '239.79', '  -1.18%  '),
    @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.001', '  +0.13%  '),
    @('LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '1.894.66', '  +0.75%  '),
    @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3031', '  -1.84%  '),
    @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '25.20', '  -4.47%  '),
    @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06812', '  -1.23%  '),
    @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07960', '  +0.20%  '),
    @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.898.48', '  +0.30%  '),
    @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.7333', '  -4.68%  '),
    @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.150', '  -1.29%  '),
    @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '90.83', '  -0.55%  '),
    @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.813.32', '  -0.54%  '),
    @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.71', '  -2.73%  '),
    @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.884', '  +1.57%  '),
    @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '241.69', '  +1.02%  '),
    @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007673', '  -1.20%  '),
    @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.04%  '),
    @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.125.84', '  -0.55%  '),
    @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  +0.10%  '),
    @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.884', '  -0.29%  '),
    @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '166.27', '  +0.75%  '),
    @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.201', '  -0.98%  '),
    @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.60', '  -0.93%  '),
    @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1278', '  +0.69%  '),
    @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.022', '  +0.10%  '),
    @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.400', '  +3.99%  '),
    @('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.512', '  -1.48%  '),
    @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.251', '  -0.92%  '),
    @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.050', '  -0.15%  '),
    @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05202', '  +2.11%  '),
    @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.246', '  -2.25%  '),
    @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7224', '  -1.50%  '),
    @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.718', '  +0.06%  '),
    @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01910', '  -0.79%  '),
    @('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.769', '  -0.52%  '),
    @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.137', '  -2.34%  '),
    @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4368', '  -1.59%  '),
    @('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '71.55', '  -3.48%  '),
    @('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.001', '  +0.10%  '),
    @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8326', '  -0.37%  '),
    @('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.873', '  -3.05%  '),
    @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.568', '  -0.98%  '),
    @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '99.55', '  -1.46%  '),
    @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.665', '  -1.36%  '),
    @('RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.038.88', '  -0.45%  '),
    @('Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '36.03', '  -0.87%  ')
)

$rowCount = $data.Count
$colCount = $data[0].Count
$values = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $values[$i, $j] = $data[$i][$j]
    }
}

# Force text formatting on the numeric-looking Price/Volume columns so
# values such as "25.20" or "1.400" keep their original digits instead of
# being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("B2:E51").Value = $values

Write-Output "Updated cryptos list."